# Refresh Leve market-board profit figures (H:N columns) across all Job sheets.
# Values below reflect the latest Universalis price snapshot pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 91485.37
$ws.Range("I2").Value = 469
$ws.Range("J2").Value = 167332.33
$ws.Range("K2").Value = 469
$ws.Range("L2").Value = 167332.33
$ws.Range("M2").Value = -356
$ws.Range("N2").Value = -167558.33
$ws.Range("H74").Value = 12500
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4064
$ws.Range("H77").Value = 12500
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20320
$ws.Range("H132").Value = 1044844.9
$ws.Range("I132").Value = 1078334.5
$ws.Range("J132").Value = 6666
$ws.Range("K132").Value = 3235003.5
$ws.Range("L132").Value = 19998
$ws.Range("M132").Value = -3232473.5
$ws.Range("N132").Value = -25058

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1707.4736
$ws.Range("I2").Value = 1825.2858
$ws.Range("J2").Value = 1377.6
$ws.Range("K2").Value = 1825.2858
$ws.Range("L2").Value = 1377.6
$ws.Range("M2").Value = -1712.2858
$ws.Range("N2").Value = -1603.6
$ws.Range("H88").Value = 1806.5454
$ws.Range("J88").Value = 1697.5714
$ws.Range("L88").Value = 1697.5714
$ws.Range("N88").Value = -2509.5714
$ws.Range("H91").Value = 1806.5454
$ws.Range("J91").Value = 1697.5714
$ws.Range("L91").Value = 1697.5714
$ws.Range("N91").Value = -4505.5714
$ws.Range("H116").Value = 1707.4736
$ws.Range("I116").Value = 1825.2858
$ws.Range("J116").Value = 1377.6
$ws.Range("K116").Value = 1825.2858
$ws.Range("L116").Value = 1377.6
$ws.Range("M116").Value = 468.7141999999999
$ws.Range("N116").Value = -5965.6
$ws.Range("H132").Value = 3206.262
$ws.Range("I132").Value = 2606.0881
$ws.Range("J132").Value = 5757
$ws.Range("K132").Value = 7818.2643
$ws.Range("L132").Value = 17271
$ws.Range("M132").Value = -5288.2643
$ws.Range("N132").Value = -22331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1707.4736
$ws.Range("I3").Value = 1825.2858
$ws.Range("J3").Value = 1377.6
$ws.Range("K3").Value = 1825.2858
$ws.Range("L3").Value = 1377.6
$ws.Range("M3").Value = -1711.2858
$ws.Range("N3").Value = -1605.6
$ws.Range("H134").Value = 3537.7144
$ws.Range("I134").Value = 2732.8215
$ws.Range("J134").Value = 6757.2856
$ws.Range("K134").Value = 8198.4645
$ws.Range("L134").Value = 20271.8568
$ws.Range("M134").Value = -5663.4645
$ws.Range("N134").Value = -25341.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 3305.4443
$ws.Range("I3").Value = 3749.5
$ws.Range("J3").Value = 3178.5715
$ws.Range("K3").Value = 3749.5
$ws.Range("L3").Value = 3178.5715
$ws.Range("M3").Value = -3636.5
$ws.Range("N3").Value = -3404.5715
$ws.Range("H4").Value = 919772.8
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1011750
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1011750
$ws.Range("M4").Value = 111
$ws.Range("N4").Value = -1011974
$ws.Range("H62").Value = 77290
$ws.Range("I62").Value = 97278.17999999999
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 97278.17999999999
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -96654.17999999999
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 77290
$ws.Range("I65").Value = 97278.17999999999
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 486390.9
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -483270.9
$ws.Range("N65").Value = -26240
$ws.Range("H132").Value = 7090.364
$ws.Range("I132").Value = 5299.4
$ws.Range("K132").Value = 15898.2
$ws.Range("M132").Value = -13368.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 30499.5
$ws.Range("I80").Value = 35666
$ws.Range("J80").Value = 15000
$ws.Range("K80").Value = 106998
$ws.Range("L80").Value = 45000
$ws.Range("M80").Value = -106062
$ws.Range("N80").Value = -46872
$ws.Range("H83").Value = 30499.5
$ws.Range("I83").Value = 35666
$ws.Range("J83").Value = 15000
$ws.Range("K83").Value = 320994
$ws.Range("L83").Value = 135000
$ws.Range("M83").Value = -316314
$ws.Range("N83").Value = -144360
$ws.Range("H116").Value = 14482.25
$ws.Range("I116").Value = 24764.5
$ws.Range("J116").Value = 4200
$ws.Range("K116").Value = 74293.5
$ws.Range("L116").Value = 12600
$ws.Range("M116").Value = -70851.5
$ws.Range("N116").Value = -19484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 12804.218
$ws.Range("I24").Value = 6083.3335
$ws.Range("J24").Value = 20136.092
$ws.Range("K24").Value = 6083.3335
$ws.Range("L24").Value = 20136.092
$ws.Range("M24").Value = -5910.3335
$ws.Range("N24").Value = -20482.092
$ws.Range("H132").Value = 1542345.2
$ws.Range("I132").Value = 1908477
$ws.Range("J132").Value = 4591.8
$ws.Range("K132").Value = 5725431
$ws.Range("L132").Value = 13775.4
$ws.Range("M132").Value = -5722901
$ws.Range("N132").Value = -18835.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3693.7646
$ws.Range("I7").Value = 3484.6428
$ws.Range("K7").Value = 3484.6428
$ws.Range("M7").Value = -3372.6428
$ws.Range("H40").Value = 4581.125
$ws.Range("I40").Value = 4378.5
$ws.Range("J40").Value = 5999.5
$ws.Range("K40").Value = 4378.5
$ws.Range("L40").Value = 5999.5
$ws.Range("M40").Value = -4242.5
$ws.Range("N40").Value = -6271.5
$ws.Range("H126").Value = 3693.7646
$ws.Range("I126").Value = 3484.6428
$ws.Range("K126").Value = 10453.9284
$ws.Range("M126").Value = -7983.928400000001
$ws.Range("H132").Value = 10264.65
$ws.Range("I132").Value = 8840.941000000001
$ws.Range("K132").Value = 26522.823
$ws.Range("M132").Value = -23992.823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 27499.5
$ws.Range("I51").Value = 34999
$ws.Range("K51").Value = 34999
$ws.Range("M51").Value = -34489
$ws.Range("H81").Value = 3969714
$ws.Range("I81").Value = 5103536
$ws.Range("K81").Value = 10207072
$ws.Range("M81").Value = -10206011
$ws.Range("H84").Value = 3969714
$ws.Range("I84").Value = 5103536
$ws.Range("K84").Value = 51035360
$ws.Range("M84").Value = -51030056
$ws.Range("H132").Value = 2779.4443
$ws.Range("I132").Value = 2730.6365
$ws.Range("J132").Value = 2856.1428
$ws.Range("K132").Value = 8191.9095
$ws.Range("L132").Value = 8568.428400000001
$ws.Range("M132").Value = -5661.9095
$ws.Range("N132").Value = -13628.4284
